$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows down: row N (for N = 11 down to 3) takes the B:F values
# currently held by row N-1 (the "rolling window" update - a new data
# point is inserted at row 2 and the oldest row's values are pushed out).
for ($r = 11; $r -ge 3; $r--) {
    $prev = $r - 1
    for ($c = 2; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($prev, $c).Value()
    }
}

# New values for row 2 (B2:F2)
$ws.Range("B2").Value = 0.01132367786385012
$ws.Range("C2").Value = 2.289151444524298
$ws.Range("D2").Value = 20.74398971997876
$ws.Range("E2").Value = 4.554557027854494
$ws.Range("F2").Value = 4.656905011860751

# Column G (count) is incremented by 1 in every data row
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 7).Value() + 1
}
